$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 48

# Column A holds a text date value - temporarily force text format so Excel
# does not auto-convert the date-like string into a date serial number, then
# restore the default cell style so no extraneous style index is recorded.
$cellA = $ws.Cells.Item($row, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "01/11/2026"
$cellA.Style = "Normal"

$ws.Cells.Item($row, 2).Value = 13216.22
$ws.Cells.Item($row, 3).Value = 0.2091078640570865
$ws.Cells.Item($row, 4).Value = 0.7908921359429135
$ws.Cells.Item($row, 5).Value = -105.47
$ws.Cells.Item($row, 6).Value = -16.61
$ws.Cells.Item($row, 7).Value = -20239.15
$ws.Cells.Item($row, 8).Value = -65.94
$ws.Cells.Item($row, 9).Value = -389.23
$ws.Cells.Item($row, 10).Value = -12.35
